$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$bg = $m.Background
Write-Host ($bg | Get-Member)
